# "add a new holy stone"
#
# GameShop.xlsx has a data table (A3:D36) listing shop items:
#   A = Id (sequential numeric id), B = Item (name), C = Shelf, D = UseDiamond
#
# A new item "fuwen-aide" (a "holy stone") is inserted right after the
# existing "fuwen-aier" row (worksheet row 10), on the same shelf (2).
# All the rows below shift down by one, and the Id column is simply kept
# sequential (it is independent bookkeeping, not tied to the item), so
# every Id from the new row down to the end increases by one, with a
# brand new Id appearing at the very end of the table.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert a new worksheet row at row 11 (just after "fuwen-aier"), pushing
# every row from 11..36 down to 12..37 - this carries the Item/Shelf/
# UseDiamond data already in those rows down with it.
[void]$ws.Rows.Item(11).Insert()

# Populate the new row with the new holy stone item.
$ws.Range("A11").Value = 15000018
$ws.Range("B11").Value = "fuwen-aide"
$ws.Range("C11").Value = 2
# Copy UseDiamond (and its style) from the row below so it matches the
# "true" value/format used throughout the rest of the table.
[void]$ws.Range("D12").Copy($ws.Range("D11"))

# The Id column (A) is just a simple running count independent of the
# item rows, so renumber it sequentially for every row pushed down by
# the insertion (this also produces a brand-new Id for the final row).
for ($r = 12; $r -le 37; $r++) {
    $ws.Cells.Item($r, 1).Value = 15000018 + ($r - 11)
}

# Grow the table (ListObject) so it covers the newly added row.
$lo = $ws.ListObjects.Item(1)
[void]$lo.Resize($ws.Range("A3:D37"))

# Leave the selection near the newly-added item, like the author did.
[void]$ws.Range("A23").Select()
